$wb = $excel.ActiveWorkbook

$oldName = "Cable Hip Abbduction"
$newName = "Cable Hip Abduction"

$ws = $wb.Worksheets.Item($oldName)

# Fix up every chart series formula on this sheet that references the old
# (misspelled) sheet name before renaming, so the stored references pick up
# the corrected name.
foreach ($co in $ws.ChartObjects()) {
    $ch = $co.Chart
    $series = $ch.SeriesCollection()
    for ($i = 1; $i -le $series.Count; $i++) {
        $s = $series.Item($i)
        $f = $s.Formula
        if ($f -like "*$oldName*") {
            $s.Formula = $f.Replace($oldName, $newName)
        }
    }
}

# Now rename the worksheet itself.
$ws.Name = $newName
